# Add a new column K "intervention_type" with a value per clinical trial row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell, matching the style of existing headers (row 1, A1:J1 use style index 1)
$ws.Range("K1").Value = "intervention_type"
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats

$values = @(
    "DEVICE",
    "PROCEDURE",
    "PROCEDURE",
    "DRUG",
    "DEVICE",
    "DIAGNOSTIC_TEST",
    "DRUG",
    "OTHER",
    "DEVICE",
    "PROCEDURE",
    "PROCEDURE",
    "PROCEDURE",
    "DRUG",
    "OTHER",
    "PROCEDURE",
    "PROCEDURE",
    "DRUG",
    "COMBINATION_PRODUCT",
    "DRUG",
    "DIAGNOSTIC_TEST"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $values[$i]
}
